# Generated edit script: applies the numeric-value updates from the
# "chore: update Sheets via scheduled runner" commit to before.xlsx.
#
# The workbook has 8 leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# each shaped identically: columns H..N are
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ.
# These are scraped/cached market-board numbers (no formulas in the sheet),
# refreshed wholesale by the scheduled runner, so the edit is just writing
# the new cached values (and, for BSM!N44, dropping the cell because the
# HQ leg no longer prices out).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 425.83334
$ws.Range("I31").Value = 425.83334
$ws.Range("K31").Value = 1277.50002
$ws.Range("M31").Value = -1047.50002

$ws.Range("H76").Value = 3892.8572
$ws.Range("J76").Value = 3892.8572
$ws.Range("L76").Value = 3892.8572
$ws.Range("N76").Value = -4522.8572

$ws.Range("H79").Value = 3892.8572
$ws.Range("J79").Value = 3892.8572
$ws.Range("L79").Value = 3892.8572
$ws.Range("N79").Value = -6076.8572

$ws.Range("H136").Value = 44595
$ws.Range("J136").Value = 44595
$ws.Range("L136").Value = 44595
$ws.Range("N136").Value = -54795

$ws.Range("H137").Value = 1663.7333
$ws.Range("I137").Value = 1386.7273
$ws.Range("J137").Value = 1928.6957
$ws.Range("K137").Value = 4160.1819
$ws.Range("L137").Value = 5786.0871
$ws.Range("M137").Value = -1610.1819
$ws.Range("N137").Value = -10886.0871


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20006690
$ws.Range("I32").Value = 34485468
$ws.Range("J32").Value = 12189
$ws.Range("K32").Value = 34485468
$ws.Range("L32").Value = 12189
$ws.Range("M32").Value = -34485181
$ws.Range("N32").Value = -12763

$ws.Range("H61").Value = 2984.25
$ws.Range("I61").Value = 2974.3333
$ws.Range("J61").Value = 3014
$ws.Range("K61").Value = 2974.3333
$ws.Range("L61").Value = 3014
$ws.Range("M61").Value = -2762.3333
$ws.Range("N61").Value = -3438

$ws.Range("H74").Value = 3076.4
$ws.Range("I74").Value = 3621.2646
$ws.Range("J74").Value = 1392.2727
$ws.Range("K74").Value = 3621.2646
$ws.Range("L74").Value = 1392.2727
$ws.Range("M74").Value = -2747.2646
$ws.Range("N74").Value = -3140.2727

$ws.Range("H77").Value = 3076.4
$ws.Range("I77").Value = 3621.2646
$ws.Range("J77").Value = 1392.2727
$ws.Range("K77").Value = 18106.323
$ws.Range("L77").Value = 6961.363499999999
$ws.Range("M77").Value = -13738.323
$ws.Range("N77").Value = -15697.3635

$ws.Range("H132").Value = 3876.762
$ws.Range("I132").Value = 3871.7144
$ws.Range("J132").Value = 3886.8572
$ws.Range("K132").Value = 11615.1432
$ws.Range("L132").Value = 11660.5716
$ws.Range("M132").Value = -9085.143199999999
$ws.Range("N132").Value = -16720.5716

$ws.Range("H136").Value = 2984.25
$ws.Range("I136").Value = 2974.3333
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 8922.999899999999
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -6372.999899999999
$ws.Range("N136").Value = -14142


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 7001
$ws.Range("I44").Value = 7001
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 7001
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H55").Value = 37000
$ws.Range("J55").Value = 37000
$ws.Range("L55").Value = 37000
$ws.Range("N55").Value = -37546

$ws.Range("H105").Value = 2998.99
$ws.Range("I105").Value = 2949.5
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2949.5
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1202.5
$ws.Range("N105").Value = -6494

$ws.Range("H134").Value = 2548.04
$ws.Range("I134").Value = 2328.9443
$ws.Range("J134").Value = 3111.4285
$ws.Range("K134").Value = 6986.8329
$ws.Range("L134").Value = 9334.2855
$ws.Range("M134").Value = -4451.8329
$ws.Range("N134").Value = -14404.2855


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4094.875
$ws.Range("I16").Value = 1760
$ws.Range("J16").Value = 7986.3335
$ws.Range("K16").Value = 1760
$ws.Range("L16").Value = 7986.3335
$ws.Range("M16").Value = -1473
$ws.Range("N16").Value = -8560.333500000001

$ws.Range("H99").Value = 2857.0952
$ws.Range("I99").Value = 2077.7778
$ws.Range("J99").Value = 3441.5833
$ws.Range("K99").Value = 2077.7778
$ws.Range("L99").Value = 3441.5833
$ws.Range("M99").Value = -579.7777999999998
$ws.Range("N99").Value = -6437.5833

$ws.Range("H113").Value = 4094.875
$ws.Range("I113").Value = 1760
$ws.Range("J113").Value = 7986.3335
$ws.Range("K113").Value = 1760
$ws.Range("L113").Value = 7986.3335
$ws.Range("M113").Value = 410
$ws.Range("N113").Value = -12326.3335

$ws.Range("H126").Value = 2857.0952
$ws.Range("I126").Value = 2077.7778
$ws.Range("J126").Value = 3441.5833
$ws.Range("K126").Value = 6233.3334
$ws.Range("L126").Value = 10324.7499
$ws.Range("M126").Value = -3763.3334
$ws.Range("N126").Value = -15264.7499

$ws.Range("H134").Value = 3038.7317
$ws.Range("I134").Value = 1791.5385
$ws.Range("J134").Value = 5200.533
$ws.Range("K134").Value = 5374.6155
$ws.Range("L134").Value = 15601.599
$ws.Range("M134").Value = -2839.6155
$ws.Range("N134").Value = -20671.599


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 487.46875
$ws.Range("I5").Value = 314.55554
$ws.Range("J5").Value = 555.13043
$ws.Range("K5").Value = 943.66662
$ws.Range("L5").Value = 1665.39129
$ws.Range("M5").Value = -831.66662
$ws.Range("N5").Value = -1889.39129

$ws.Range("H132").Value = 722258.6
$ws.Range("I132").Value = 479.66666
$ws.Range("J132").Value = 1263592.9
$ws.Range("K132").Value = 4316.99994
$ws.Range("L132").Value = 11372336.1
$ws.Range("M132").Value = -1786.99994
$ws.Range("N132").Value = -11377396.1

$ws.Range("H135").Value = 487.46875
$ws.Range("I135").Value = 314.55554
$ws.Range("J135").Value = 555.13043
$ws.Range("K135").Value = 2830.99986
$ws.Range("L135").Value = 4996.173870000001
$ws.Range("M135").Value = -295.9998599999999
$ws.Range("N135").Value = -10066.17387

$ws.Range("H136").Value = 3291.25
$ws.Range("I136").Value = 1666
$ws.Range("K136").Value = 4998
$ws.Range("M136").Value = 102


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 15720.2
$ws.Range("I41").Value = 2136
$ws.Range("K41").Value = 2136
$ws.Range("M41").Value = -1781

$ws.Range("H132").Value = 5688.875
$ws.Range("I132").Value = 6001.769
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 18005.307
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -15475.307
$ws.Range("N132").Value = -18059


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("N50").Value = -11274

$ws.Range("H136").Value = 1582.2916
$ws.Range("I136").Value = 1028.5714
$ws.Range("J136").Value = 2357.5
$ws.Range("K136").Value = 3085.7142
$ws.Range("L136").Value = 7072.5
$ws.Range("M136").Value = -535.7142000000003
$ws.Range("N136").Value = -12172.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 18321.4
$ws.Range("J9").Value = 18321.4
$ws.Range("L9").Value = 18321.4
$ws.Range("N9").Value = -18601.4

$ws.Range("H122").Value = 2422.138
$ws.Range("I122").Value = 1641.4117
$ws.Range("J122").Value = 3528.1667
$ws.Range("K122").Value = 4924.2351
$ws.Range("L122").Value = 10584.5001
$ws.Range("M122").Value = -2474.2351
$ws.Range("N122").Value = -15484.5001

$ws.Range("H126").Value = 1107.9565
$ws.Range("I126").Value = 1076.3125
$ws.Range("J126").Value = 1180.2858
$ws.Range("K126").Value = 3228.9375
$ws.Range("L126").Value = 3540.8574
$ws.Range("M126").Value = -758.9375
$ws.Range("N126").Value = -8480.857400000001

$ws.Range("H132").Value = 2109.8286
$ws.Range("I132").Value = 1458.7858
$ws.Range("J132").Value = 4714
$ws.Range("K132").Value = 4376.357400000001
$ws.Range("L132").Value = 14142
$ws.Range("M132").Value = -1846.357400000001
$ws.Range("N132").Value = -19202

$ws.Range("H136").Value = 40577.92
$ws.Range("I136").Value = 63938.5
$ws.Range("J136").Value = 3201
$ws.Range("K136").Value = 191815.5
$ws.Range("L136").Value = 9603
$ws.Range("M136").Value = -189265.5
$ws.Range("N136").Value = -14703
